# Apply the "quote the literal operand" edit to the String functions example
# column (C44:C59) on Sheet1, then leave the view scrolled/selected the way
# the author last left it (top-left C37, active cell C59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    44 = '"helloworld".contains("hello")'
    45 = '"hello world".endsWith("old")'
    46 = '"foo".matches("k.*")'
    47 = '"hello".toUpper()'
    48 = '"HEllO".toLower()'
    49 = '"hello".base64Encode()'
    50 = '"aGVsbG8=".base64Decode()'
    51 = '"hello_world!".base64URLEncode()'
    52 = '"aGVsbG9fd29ybGQh".base64URLDecode()'
    53 = '"hello".size()'
    54 = '"hello".substring(1,4)'
    55 = '"hello".split("e")'
    56 = '"hello".replaceAll("l","p")'
    57 = '"some12#$text".matchAndReplaceAll("[^a-zA-Z]+", "-")'
    58 = '"hello".indexOf("l")'
    59 = '"hello".lastIndexOf("l")'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Match the saved view state from the diff: scrolled to C37, active cell C59.
$ws.Application.Goto($ws.Range("C59"), $true)
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 3
